$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Kevytkoris" label next to the existing 10.00-12.00 Sunday slot
# for Maininki (row 8, column O = Sunday second column)
$ws.Range("O8").Value = "Kevytkoris"
